$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.093010187149048
$ws.Range("B1").Value = 2.174811124801636
$ws.Range("C1").Value = 9.596397399902344
$ws.Range("D1").Value = 1.123750448226929
$ws.Range("E1").Value = 1.210295438766479
